$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '78.985.21'
$ws.Range("E2").Value = '  +3.30%  '

$ws.Range("D3").Value = '3.185.86'

$ws.Range("E4").Value = '  -0.05%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '205.10'
$ws.Range("E5").Value = '  +2.50%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '630.15'
$ws.Range("E6").Value = '  +0.46%  '

$ws.Range("E7").Value = '  +0.01%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.228'
$ws.Range("E8").Value = '  +9.44%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.586'
$ws.Range("E9").Value = '  +6.12%  '

$ws.Range("D10").Value = '3.186.27'
$ws.Range("E10").Value = '  +5.35%  '

$ws.Range("E11").Value = '  +34.64%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.165'
$ws.Range("E12").Value = '  +2.77%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.45'
$ws.Range("E13").Value = '  +5.81%  '

$ws.Range("D14").Value = '3.775.07'
$ws.Range("E14").Value = '  +5.32%  '

$ws.Range("E15").Value = '  +17.54%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '31.60'
$ws.Range("E16").Value = '  +7.95%  '

$ws.Range("D17").Value = '78.856.01'
$ws.Range("E17").Value = '  +3.30%  '

$ws.Range("D18").Value = '3.183.21'
$ws.Range("E18").Value = '  +5.17%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '14.50'
$ws.Range("E19").Value = '  +7.69%  '

$ws.Range("E20").Value = '  +2.93%  '

$ws.Range("B21").Value = 'SuiNetwork'
$ws.Range("C21").Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '2.89'
$ws.Range("E21").Value = '  +27.37%  '

$ws.Range("B22").Value = 'BitcoinCash'
$ws.Range("C22").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '427.96'
$ws.Range("E22").Value = '  +14.61%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.97'
$ws.Range("E23").Value = '  +14.24%  '

$ws.Range("E24").Value = '  +5.95%  '

$ws.Range("D25").Value = '3.351.52'
$ws.Range("E25").Value = '  +5.36%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '4.78'
$ws.Range("E26").Value = '  +9.16%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '11.08'
$ws.Range("E27").Value = '  +11.95%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '76.06'
$ws.Range("E28").Value = '  +3.99%  '

$ws.Range("E29").Value = '  +0.26%  '

$ws.Range("E30").Value = '  +3.76%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.995'
$ws.Range("E31").Value = '  -0.31%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '8.92'
$ws.Range("E32").Value = '  +7.48%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.48'
$ws.Range("E33").Value = '  +4.86%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '520.22'
$ws.Range("E34").Value = '  +2.71%  '

$ws.Range("E35").Value = '  +2.29%  '

$ws.Range("B36").Value = 'EthereumClassic'
$ws.Range("C36").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '23.00'
$ws.Range("E36").Value = '  +11.11%  '

$ws.Range("B37").Value = 'Cronos'
$ws.Range("C37").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.127'
$ws.Range("E37").Value = '  +20.66%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.135'
$ws.Range("E38").Value = '  +20.13%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.399'
$ws.Range("E40").Value = '  +3.83%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '164.06'
$ws.Range("E41").Value = '  +0.05%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '19.99'
$ws.Range("E42").Value = '  -0.04%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '192.64'
$ws.Range("E43").Value = '  +1.35%  '

$ws.Range("E44").Value = '  +0.01%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '5.43'
$ws.Range("E45").Value = '  +5.92%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.808'
$ws.Range("E46").Value = '  +13.13%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.78'
$ws.Range("E47").Value = '  +6.66%  '

$ws.Range("E48").Value = '  +4.28%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '42.78'
$ws.Range("E49").Value = '  +0.88%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.50'
$ws.Range("E50").Value = '  +5.60%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.624'
$ws.Range("E51").Value = '  +3.64%  '
